$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("H 72") and shift all rows below it up by one.
$ws.Rows("2:2").Delete()
